$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (bold) ----
$ws.Range("B4").Value = "Account Number"
$ws.Range("C4").Value = "Amount"
$ws.Range("B4:C4").Font.Bold = $true

# ---- Account number / amount table ----
# The first two account numbers were originally typed with a leading
# apostrophe (quote-prefix => forces text), the remaining ones were typed
# as plain text. All account numbers end up right aligned.
$ws.Range("B5").Value = "'A0000"
$ws.Range("B5").HorizontalAlignment = -4152
$ws.Range("C5").Value = 4747

$ws.Range("B6").HorizontalAlignment = -4152
$ws.Range("B6").Value = "'A1111"
$ws.Range("C6").Value = 0

$ws.Range("B7").HorizontalAlignment = -4152
$ws.Range("B7").Value = "A2222"
$ws.Range("C7").Value = 5228

$ws.Range("B8").HorizontalAlignment = -4152
$ws.Range("B8").Value = "A3333"
$ws.Range("C8").Value = 0

$ws.Range("B9").HorizontalAlignment = -4152
$ws.Range("B9").Value = "A4444"
$ws.Range("C9").Value = 4741

$ws.Range("B10").HorizontalAlignment = -4152
$ws.Range("B10").Value = "A5555"
$ws.Range("C10").Value = 4445

$ws.Range("B11").HorizontalAlignment = -4152
$ws.Range("B11").Value = "A7777"
$ws.Range("C11").Value = 9560

$ws.Range("B12").HorizontalAlignment = -4152
$ws.Range("B12").Value = "A8888"
$ws.Range("C12").Value = 0

$ws.Range("B13").HorizontalAlignment = -4152
$ws.Range("B13").Value = "A9999"
$ws.Range("C13").Value = 0

# ---- Column widths (best-fit to contents) ----
$ws.Range("B1").ColumnWidth = 15.166666666666666
$ws.Range("C1").ColumnWidth = 15.166666666666666
$ws.Range("D1").ColumnWidth = 7.333333333333333

# ---- Selection, matching where the author left the cursor ----
$ws.Range("H11").Select() | Out-Null
